$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out the old "reference / sanity-check" block (rows 11, 13-15) ---
$ws.Range("A11:XFD15").ClearContents()

# --- Add F column conversion formulas (acres/10000 -> km-ish) for data rows 2-7 ---
$ws.Range("F2").Formula = "=D2/10000"
$ws.Range("F3:F7").Formula = "=D3/10000"

# --- Row 8: sum of F2:F7 ---
$ws.Range("F8").Formula = "=SUM(F2:F7)"

# --- Row 9: keep existing A9 text, add F9 conversion to km ---
$ws.Range("F9").Formula = "=F8/1000"

# --- New M column: header + USGS km values for a few regions ---
$ws.Range("M1").Value = "km_USGS"
$ws.Range("M2").Value = 90
$ws.Range("M5").Value = 107
$ws.Range("M6").Value = 75

# --- Update the selected cell / window view to match the author's last position ---
$ws.Range("M7").Select()
$wb.Windows.Item(1).Left = -165
$wb.Windows.Item(1).Top = 6285

$wb.Save()
